# This workbook logs daily "Poroto verde" price observations, one row per
# observation, ordered (mostly) from newest to oldest. Two new observations
# were recorded and inserted at the top of the data block (rows 53-54,
# right after the existing row 52), pushing the previously-existing rows
# 53-147 down to rows 55-149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 53:54 - this shifts old rows 53-147 down to
# 55-149 and grows the sheet dimension from R147 to R149 automatically.
$ws.Range("A53:A54").EntireRow.Insert()

# New row 53
$ws.Cells.Item(53, 1).Value = 8
$ws.Cells.Item(53, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(53, 3).Value = "Coquimbo"
$ws.Cells.Item(53, 4).Value = 44533
$ws.Cells.Item(53, 5).Value = 4
$ws.Cells.Item(53, 6).Value = 100112031
$ws.Cells.Item(53, 7).Value = "Poroto verde"
$ws.Cells.Item(53, 8).Value = "Magnum"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 540
$ws.Cells.Item(53, 11).Value = 17000
$ws.Cells.Item(53, 12).Value = 18000
$ws.Cells.Item(53, 13).Value = 17500
$ws.Cells.Item(53, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(53, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(53, 16).Value = 700
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 18).Value = "Hortaliza"

# New row 54
$ws.Cells.Item(54, 1).Value = 8
$ws.Cells.Item(54, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value = 44533
$ws.Cells.Item(54, 5).Value = 4
$ws.Cells.Item(54, 6).Value = 100112031
$ws.Cells.Item(54, 7).Value = "Poroto verde"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 400
$ws.Cells.Item(54, 11).Value = 27000
$ws.Cells.Item(54, 12).Value = 28000
$ws.Cells.Item(54, 13).Value = 27500
$ws.Cells.Item(54, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(54, 16).Value = 1100
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = "Hortaliza"

# Column D (Fecha) carries a date display style (s="2" in the original
# workbook); make sure the two freshly-inserted rows keep/display the same
# number format as the surrounding date cells.
$ws.Cells.Item(53, 4).NumberFormat = $ws.Cells.Item(55, 4).NumberFormat
$ws.Cells.Item(54, 4).NumberFormat = $ws.Cells.Item(55, 4).NumberFormat
